# Nexial unitTest_function.xlsx update:
#  - rdbms: add assertResultMatch(var,columns,search) and assertResultNotMatch(var,columns,search)
#  - localdb: add queryAsCSV(var,sql)
#  - web: add assertElementDisabled(locator); rename checkAll/uncheckAll to add waitMs param
#  - step.inTime: new category (observe/perform/validate) reusing column Z, replacing old tn.5250 data
#  - defined names: update localdb/rdbms/web ranges, add step.inTime

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------
# 1) localdb (column O): insert "queryAsCSV(var,sql)" before "runSQLs(var,sqls)"
#    O2:O12 (11 items) -> O2:O13 (12 items)
# ---------------------------------------------------------------
$ws.Range("O13").Value = $ws.Range("O12").Value2
$ws.Range("O12").Value = "queryAsCSV(var,sql)"

# ---------------------------------------------------------------
# 2) rdbms (column T): insert 2 new commands at the top (alphabetically first)
#    T2:T7 (6 items) -> T2:T9 (8 items)
# ---------------------------------------------------------------
$ws.Range("T9").Value = $ws.Range("T7").Value2
$ws.Range("T8").Value = $ws.Range("T6").Value2
$ws.Range("T7").Value = $ws.Range("T5").Value2
$ws.Range("T6").Value = $ws.Range("T4").Value2
$ws.Range("T5").Value = $ws.Range("T3").Value2
$ws.Range("T4").Value = $ws.Range("T2").Value2
$ws.Range("T3").Value = "assertResultNotMatch(var,columns,search)"
$ws.Range("T2").Value = "assertResultMatch(var,columns,search)"

# ---------------------------------------------------------------
# 3) column Z: replace "tn.5250" (close/open/saveText/typeKeys/updateScreenFields)
#    content with new "step.inTime" (observe/perform/validate) content, and
#    clear the two now-unused trailing cells.
# ---------------------------------------------------------------
$ws.Range("Z1").Value = "step.inTime"
$ws.Range("Z2").Value = "observe(prompt,waitMs)"
$ws.Range("Z3").Value = "perform(instructions,waitMs)"
$ws.Range("Z4").Value = "validate(prompt,responses,passResponses,waitMs)"
$ws.Range("Z5").ClearContents()
$ws.Range("Z6").ClearContents()

# ---------------------------------------------------------------
# 4) web (column AA): insert "assertElementDisabled(locator)" before
#    "assertElementEnabled(locator)" (row 15), shifting everything below
#    down by one row: AA15:AA144 (130 items) -> AA16:AA145 (131 items)
# ---------------------------------------------------------------
for ($i = 144; $i -ge 15; $i--) {
    $src = $ws.Range("AA$i").Value2
    $dst = $i + 1
    $ws.Range("AA$dst").Value = $src
}
$ws.Range("AA15").Value = "assertElementDisabled(locator)"

# rename checkAll(locator) -> checkAll(locator,waitMs) (same alphabetical slot, now row 51)
$ws.Range("AA51").Value = "checkAll(locator,waitMs)"

# rename uncheckAll(locator) -> uncheckAll(locator,waitMs) (same alphabetical slot, now row 134)
$ws.Range("AA134").Value = "uncheckAll(locator,waitMs)"

# ---------------------------------------------------------------
# 5) defined names: update ranges that grew, add the new step.inTime name
# ---------------------------------------------------------------
$wb.Names.Item("localdb").RefersTo = "='#system'!`$O`$2:`$O`$13"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$145"
$wb.Names.Add("step.inTime", "='#system'!`$Z`$2:`$Z`$4")
